# Auto-generated cell updates for cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Delete()
}

$ws.Range("D2").Value = '29.402.50'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.841.88'
Set-TextValue $ws.Range("D4") '0.9988'
$ws.Range("E4").Value = '  +0.09%  '
Set-TextValue $ws.Range("D5") '239.65'
$ws.Range("E5").Value = '  -0.35%  '
Set-TextValue $ws.Range("D6") '0.6270'
$ws.Range("E6").Value = '  -0.04%  '
Set-TextValue $ws.Range("D8") '0.07447'
$ws.Range("E8").Value = '  -0.69%  '
Set-TextValue $ws.Range("D9") '0.2901'
$ws.Range("E9").Value = '  -0.20%  '
Set-TextValue $ws.Range("D10") '24.81'
$ws.Range("E10").Value = '  +1.45%  '
Set-TextValue $ws.Range("D11") '0.07718'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '1.835.54'
$ws.Range("E12").Value = '  -0.63%  '
Set-TextValue $ws.Range("D13") '4.970'
$ws.Range("E13").Value = '  -0.63%  '
Set-TextValue $ws.Range("D14") '0.6767'
$ws.Range("E14").Value = '  -0.64%  '
Set-TextValue $ws.Range("D15") '0.00001025'
$ws.Range("E15").Value = '  -3.05%  '
Set-TextValue $ws.Range("D16") '81.82'
$ws.Range("E16").Value = '  -0.50%  '
Set-TextValue $ws.Range("D17") '6.250'
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").Value = '29.451.57'
$ws.Range("E18").Value = '  +0.08%  '
Set-TextValue $ws.Range("D19") '232.11'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("E21").Value = '  +0.12%  '
Set-TextValue $ws.Range("D22") '7.326'
$ws.Range("E22").Value = '  -2.17%  '
Set-TextValue $ws.Range("D23") '1.001'
$ws.Range("E23").Value = '  +0.11%  '
Set-TextValue $ws.Range("D24") '158.19'
$ws.Range("E24").Value = '  -0.59%  '
Set-TextValue $ws.Range("D25") '8.489'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("E26").Value = '  -1.82%  '
Set-TextValue $ws.Range("D27") '17.37'
$ws.Range("E27").Value = '  -1.00%  '
Set-TextValue $ws.Range("D28") '0.07280'
$ws.Range("E28").Value = '  +12.59%  '
Set-TextValue $ws.Range("D29") '1.461'
$ws.Range("E29").Value = '  +3.05%  '
Set-TextValue $ws.Range("D30") '1.480'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  -1.17%  '
Set-TextValue $ws.Range("D32") '4.056'
$ws.Range("E32").Value = '  -1.06%  '
Set-TextValue $ws.Range("D33") '1.819'
$ws.Range("E33").Value = '  -0.66%  '
Set-TextValue $ws.Range("D34") '1.141'
$ws.Range("E34").Value = '  -0.20%  '
Set-TextValue $ws.Range("D35") '0.6976'
$ws.Range("E35").Value = '  +0.12%  '
Set-TextValue $ws.Range("D36") '2.568'
$ws.Range("E36").Value = '  -0.44%  '
Set-TextValue $ws.Range("D37") '6.983'
$ws.Range("E37").Value = '  +3.34%  '
$ws.Range("E38").Value = '  +0.20%  '
Set-TextValue $ws.Range("D39") '2.816'
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").Value = '1.235.44'
$ws.Range("E40").Value = '  -2.90%  '
Set-TextValue $ws.Range("D41") '0.9452'
Set-TextValue $ws.Range("D42") '0.9999'
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = '2.001.48'
$ws.Range("E43").Value = '  -0.44%  '
Set-TextValue $ws.Range("D44") '100.86'
Set-TextValue $ws.Range("D45") '65.61'
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("E46").Value = '  +0.33%  '
Set-TextValue $ws.Range("D47") '1.726'
$ws.Range("E47").Value = '  -1.14%  '
Set-TextValue $ws.Range("D48") '6.963'
$ws.Range("E48").Value = '  -1.73%  '
Set-TextValue $ws.Range("D49") '8.925'
$ws.Range("E49").Value = '  -1.65%  '
Set-TextValue $ws.Range("D50") '0.1139'
$ws.Range("E50").Value = '  -3.08%  '
Set-TextValue $ws.Range("D51") '0.3904'
$ws.Range("E51").Value = '  -1.33%  '
